# Atualiza todos os resources para "active"
# The Metadata sheet keeps a Property/Value table; the "Status" row currently
# holds "draft" and must be updated to "active".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")
$ws.Range("B6").Value = "active"
